$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1058.3158  # H40
$ws.Cells.Item(40, 9).Value = 1058.3158  # I40
$ws.Cells.Item(40, 10).Value = 0  # J40
$ws.Cells.Item(40, 11).Value = 1058.3158  # K40
$ws.Cells.Item(40, 12).Value = 0  # L40
$ws.Cells.Item(40, 13).Value = -883.3158000000001  # M40
$ws.Cells.Item(40, 14).ClearContents()  # N40
$ws.Cells.Item(76, 8).Value = 3333  # H76
$ws.Cells.Item(76, 10).Value = 3333  # J76
$ws.Cells.Item(76, 12).Value = 3333  # L76
$ws.Cells.Item(76, 14).Value = -3963  # N76
$ws.Cells.Item(79, 8).Value = 3333  # H79
$ws.Cells.Item(79, 10).Value = 3333  # J79
$ws.Cells.Item(79, 12).Value = 3333  # L79
$ws.Cells.Item(79, 14).Value = -5517  # N79
$ws.Cells.Item(86, 8).Value = 1081.9286  # H86
$ws.Cells.Item(86, 9).Value = 1028.2858  # I86
$ws.Cells.Item(86, 10).Value = 1135.5714  # J86
$ws.Cells.Item(86, 11).Value = 1028.2858  # K86
$ws.Cells.Item(86, 12).Value = 1135.5714  # L86
$ws.Cells.Item(86, 13).Value = 94.71419999999989  # M86
$ws.Cells.Item(86, 14).Value = -3381.5714  # N86
$ws.Cells.Item(89, 8).Value = 1081.9286  # H89
$ws.Cells.Item(89, 9).Value = 1028.2858  # I89
$ws.Cells.Item(89, 10).Value = 1135.5714  # J89
$ws.Cells.Item(89, 11).Value = 5141.429  # K89
$ws.Cells.Item(89, 12).Value = 5677.857  # L89
$ws.Cells.Item(89, 13).Value = 474.5709999999999  # M89
$ws.Cells.Item(89, 14).Value = -16909.857  # N89
$ws.Cells.Item(96, 8).Value = 1453.125  # H96
$ws.Cells.Item(96, 9).Value = 1819.5  # I96
$ws.Cells.Item(96, 11).Value = 5458.5  # K96
$ws.Cells.Item(96, 13).Value = -4085.5  # M96
$ws.Cells.Item(98, 8).Value = 2399.8438  # H98
$ws.Cells.Item(98, 9).Value = 2375.724  # I98
$ws.Cells.Item(98, 11).Value = 2375.724  # K98
$ws.Cells.Item(98, 13).Value = -877.7240000000002  # M98
$ws.Cells.Item(99, 8).Value = 1064.5883  # H99
$ws.Cells.Item(99, 9).Value = 467.18182  # I99
$ws.Cells.Item(99, 10).Value = 2159.8333  # J99
$ws.Cells.Item(99, 11).Value = 1401.54546  # K99
$ws.Cells.Item(99, 12).Value = 6479.499899999999  # L99
$ws.Cells.Item(99, 13).Value = 96.45453999999995  # M99
$ws.Cells.Item(99, 14).Value = -9475.499899999999  # N99
$ws.Cells.Item(122, 8).Value = 2399.8438  # H122
$ws.Cells.Item(122, 9).Value = 2375.724  # I122
$ws.Cells.Item(122, 11).Value = 7127.172  # K122
$ws.Cells.Item(122, 13).Value = -4677.172  # M122
$ws.Cells.Item(138, 8).Value = 1830.1094  # H138
$ws.Cells.Item(138, 9).Value = 1570.1765  # I138
$ws.Cells.Item(138, 10).Value = 2849.8462  # J138
$ws.Cells.Item(138, 11).Value = 4710.529500000001  # K138
$ws.Cells.Item(138, 12).Value = 8549.5386  # L138
$ws.Cells.Item(138, 13).Value = 429.4704999999994  # M138
$ws.Cells.Item(138, 14).Value = -18829.5386  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3055.3452  # H32
$ws.Cells.Item(32, 9).Value = 2570.7375  # I32
$ws.Cells.Item(32, 11).Value = 2570.7375  # K32
$ws.Cells.Item(32, 13).Value = -2283.7375  # M32
$ws.Cells.Item(74, 8).Value = 1523.3438  # H74
$ws.Cells.Item(74, 9).Value = 1250  # I74
$ws.Cells.Item(74, 11).Value = 1250  # K74
$ws.Cells.Item(74, 13).Value = -376  # M74
$ws.Cells.Item(77, 8).Value = 1523.3438  # H77
$ws.Cells.Item(77, 9).Value = 1250  # I77
$ws.Cells.Item(77, 11).Value = 6250  # K77
$ws.Cells.Item(77, 13).Value = -1882  # M77
$ws.Cells.Item(122, 8).Value = 1544.5667  # H122
$ws.Cells.Item(122, 9).Value = 1417.7273  # I122
$ws.Cells.Item(122, 11).Value = 4253.1819  # K122
$ws.Cells.Item(122, 13).Value = -1803.1819  # M122

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3433.2444  # H134
$ws.Cells.Item(134, 9).Value = 3355.2092  # I134
$ws.Cells.Item(134, 11).Value = 10065.6276  # K134
$ws.Cells.Item(134, 13).Value = -7530.6276  # M134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2557.4375  # H31
$ws.Cells.Item(31, 9).Value = 2374.7144  # I31
$ws.Cells.Item(31, 10).Value = 2699.5557  # J31
$ws.Cells.Item(31, 11).Value = 2374.7144  # K31
$ws.Cells.Item(31, 12).Value = 2699.5557  # L31
$ws.Cells.Item(31, 13).Value = -2079.7144  # M31
$ws.Cells.Item(31, 14).Value = -3289.5557  # N31
$ws.Cells.Item(34, 8).Value = 2557.4375  # H34
$ws.Cells.Item(34, 9).Value = 2374.7144  # I34
$ws.Cells.Item(34, 10).Value = 2699.5557  # J34
$ws.Cells.Item(34, 11).Value = 2374.7144  # K34
$ws.Cells.Item(34, 12).Value = 2699.5557  # L34
$ws.Cells.Item(34, 13).Value = -2172.7144  # M34
$ws.Cells.Item(34, 14).Value = -3103.5557  # N34
$ws.Cells.Item(58, 8).Value = 853510.2  # H58
$ws.Cells.Item(58, 9).Value = 1175715.2  # I58
$ws.Cells.Item(58, 11).Value = 1175715.2  # K58
$ws.Cells.Item(58, 13).Value = -1175512.2  # M58
$ws.Cells.Item(62, 8).Value = 3299.889  # H62
$ws.Cells.Item(62, 9).Value = 3160  # I62
$ws.Cells.Item(62, 10).Value = 3474.75  # J62
$ws.Cells.Item(62, 11).Value = 3160  # K62
$ws.Cells.Item(62, 12).Value = 3474.75  # L62
$ws.Cells.Item(62, 13).Value = -2536  # M62
$ws.Cells.Item(62, 14).Value = -4722.75  # N62
$ws.Cells.Item(65, 8).Value = 3299.889  # H65
$ws.Cells.Item(65, 9).Value = 3160  # I65
$ws.Cells.Item(65, 10).Value = 3474.75  # J65
$ws.Cells.Item(65, 11).Value = 15800  # K65
$ws.Cells.Item(65, 12).Value = 17373.75  # L65
$ws.Cells.Item(65, 13).Value = -12680  # M65
$ws.Cells.Item(65, 14).Value = -23613.75  # N65
$ws.Cells.Item(132, 8).Value = 1331.3276  # H132
$ws.Cells.Item(132, 9).Value = 1088.6222  # I132
$ws.Cells.Item(132, 11).Value = 3265.8666  # K132
$ws.Cells.Item(132, 13).Value = -735.8666000000003  # M132
$ws.Cells.Item(136, 8).Value = 853510.2  # H136
$ws.Cells.Item(136, 9).Value = 1175715.2  # I136
$ws.Cells.Item(136, 11).Value = 3527145.6  # K136
$ws.Cells.Item(136, 13).Value = -3524595.6  # M136

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 326806.84  # H4
$ws.Cells.Item(4, 9).Value = 62615  # I4
$ws.Cells.Item(4, 11).Value = 187845  # K4
$ws.Cells.Item(4, 13).Value = -187733  # M4
$ws.Cells.Item(75, 8).Value = 20000  # H75
$ws.Cells.Item(75, 10).Value = 20000  # J75
$ws.Cells.Item(75, 12).Value = 60000  # L75
$ws.Cells.Item(75, 14).Value = -61996  # N75
$ws.Cells.Item(78, 8).Value = 20000  # H78
$ws.Cells.Item(78, 10).Value = 20000  # J78
$ws.Cells.Item(78, 12).Value = 180000  # L78
$ws.Cells.Item(78, 14).Value = -189984  # N78
$ws.Cells.Item(103, 8).Value = 2185.4666  # H103
$ws.Cells.Item(103, 10).Value = 1893  # J103
$ws.Cells.Item(103, 12).Value = 5679  # L103
$ws.Cells.Item(103, 14).Value = -7437  # N103
$ws.Cells.Item(130, 8).Value = 1882  # H130
$ws.Cells.Item(130, 9).Value = 1342.6666  # I130
$ws.Cells.Item(130, 11).Value = 4027.9998  # K130
$ws.Cells.Item(130, 13).Value = 992.0001999999999  # M130
$ws.Cells.Item(140, 8).Value = 2191.3518  # H140
$ws.Cells.Item(140, 9).Value = 1149.2  # I140
$ws.Cells.Item(140, 10).Value = 2592.1794  # J140
$ws.Cells.Item(140, 11).Value = 3447.6  # K140
$ws.Cells.Item(140, 12).Value = 7776.5382  # L140
$ws.Cells.Item(140, 13).Value = 1732.4  # M140
$ws.Cells.Item(140, 14).Value = -18136.5382  # N140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 29999  # H57
$ws.Cells.Item(57, 10).Value = 29999  # J57
$ws.Cells.Item(57, 12).Value = 29999  # L57
$ws.Cells.Item(57, 14).Value = -31639  # N57
$ws.Cells.Item(107, 8).Value = 100  # H107
$ws.Cells.Item(107, 9).Value = 100  # I107
$ws.Cells.Item(107, 11).Value = 100  # K107
$ws.Cells.Item(107, 13).Value = 1820  # M107
$ws.Cells.Item(113, 8).Value = 1882  # H113
$ws.Cells.Item(113, 9).Value = 2136.6667  # I113
$ws.Cells.Item(113, 11).Value = 2136.6667  # K113
$ws.Cells.Item(113, 13).Value = 33.33329999999978  # M113

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(96, 8).Value = 85000  # H96
$ws.Cells.Item(96, 10).Value = 85000  # J96
$ws.Cells.Item(96, 12).Value = 85000  # L96
$ws.Cells.Item(96, 14).Value = -90492  # N96
$ws.Cells.Item(122, 8).Value = 2833.4119  # H122
$ws.Cells.Item(122, 9).Value = 1766.3  # I122
$ws.Cells.Item(122, 10).Value = 4357.857  # J122
$ws.Cells.Item(122, 11).Value = 5298.9  # K122
$ws.Cells.Item(122, 12).Value = 13073.571  # L122
$ws.Cells.Item(122, 13).Value = -2848.9  # M122
$ws.Cells.Item(122, 14).Value = -17973.571  # N122
$ws.Cells.Item(132, 8).Value = 1559.5763  # H132
$ws.Cells.Item(132, 9).Value = 1285.1471  # I132
$ws.Cells.Item(132, 10).Value = 1932.8  # J132
$ws.Cells.Item(132, 11).Value = 3855.4413  # K132
$ws.Cells.Item(132, 12).Value = 5798.4  # L132
$ws.Cells.Item(132, 13).Value = -1325.4413  # M132
$ws.Cells.Item(132, 14).Value = -10858.4  # N132

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3191.5  # H62
$ws.Cells.Item(62, 9).Value = 2829  # I62
$ws.Cells.Item(62, 11).Value = 2829  # K62
$ws.Cells.Item(62, 13).Value = -2205  # M62
$ws.Cells.Item(65, 8).Value = 3191.5  # H65
$ws.Cells.Item(65, 9).Value = 2829  # I65
$ws.Cells.Item(65, 11).Value = 14145  # K65
$ws.Cells.Item(65, 13).Value = -11025  # M65
$ws.Cells.Item(132, 8).Value = 1105.6271  # H132
$ws.Cells.Item(132, 10).Value = 1929.7333  # J132
$ws.Cells.Item(132, 12).Value = 5789.199900000001  # L132
$ws.Cells.Item(132, 14).Value = -10849.1999  # N132

Write-Output "Applied all Tonberry_Profits value updates"